$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("H67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("H93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("H132").Value = 1750
$ws.Range("K132").Value = 2250
$ws.Range("M132").Value = 280
$ws.Range("I132").Value = 750

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 322.0909
$ws.Range("K32").Value = 322.0909
$ws.Range("M32").Value = -35.09089999999998
$ws.Range("I32").Value = 322.0909
$ws.Range("H61").Value = 2250
$ws.Range("K61").Value = 2250
$ws.Range("M61").Value = -2038
$ws.Range("I61").Value = 2250
$ws.Range("H74").Value = 3589.3
$ws.Range("N74").ClearContents()
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("H77").Value = 3589.3
$ws.Range("N77").ClearContents()
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("H92").Value = 29999.5
$ws.Range("N92").Value = -34991.5
$ws.Range("J92").Value = 29999.5
$ws.Range("L92").Value = 29999.5
$ws.Range("H102").Value = 3502.6667
$ws.Range("K102").Value = 3502.6667
$ws.Range("M102").Value = -1880.6667
$ws.Range("I102").Value = 3502.6667
$ws.Range("H132").Value = 2636.75
$ws.Range("K132").Value = 7910.25
$ws.Range("M132").Value = -5380.25
$ws.Range("I132").Value = 2636.75
$ws.Range("H135").Value = 5027247.5
$ws.Range("N135").Value = -5037387.5
$ws.Range("J135").Value = 5027247.5
$ws.Range("L135").Value = 5027247.5
$ws.Range("H136").Value = 2250
$ws.Range("K136").Value = 6750
$ws.Range("M136").Value = -4200
$ws.Range("I136").Value = 2250

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3133.3635
$ws.Range("N20").Value = -3360.8572
$ws.Range("J20").Value = 2866.8572
$ws.Range("L20").Value = 2866.8572
$ws.Range("H92").Value = 59999
$ws.Range("N92").Value = -64991
$ws.Range("J92").Value = 59999
$ws.Range("L92").Value = 59999
$ws.Range("H94").Value = 2365.8235
$ws.Range("N94").Value = -3725.3333
$ws.Range("K94").Value = 1851.125
$ws.Range("J94").Value = 2823.3333
$ws.Range("M94").Value = -1400.125
$ws.Range("L94").Value = 2823.3333
$ws.Range("I94").Value = 1851.125
$ws.Range("H134").Value = 3720
$ws.Range("K134").Value = 11160
$ws.Range("M134").Value = -8625
$ws.Range("I134").Value = 3720

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 26250000
$ws.Range("K6").Value = 26250000
$ws.Range("M6").Value = -26249887
$ws.Range("I6").Value = 26250000
$ws.Range("H7").Value = 114.57143
$ws.Range("N7").Value = -448.5
$ws.Range("K7").Value = 71.40000000000001
$ws.Range("J7").Value = 222.5
$ws.Range("M7").Value = 41.59999999999999
$ws.Range("L7").Value = 222.5
$ws.Range("I7").Value = 71.40000000000001
$ws.Range("H17").Value = 3500
$ws.Range("N17").Value = -3848
$ws.Range("K17").Value = 0
$ws.Range("J17").Value = 3500
$ws.Range("M17").ClearContents()
$ws.Range("L17").Value = 3500
$ws.Range("I17").Value = 0
$ws.Range("H22").Value = 1850
$ws.Range("N22").Value = -2550
$ws.Range("J22").Value = 1850
$ws.Range("L22").Value = 1850
$ws.Range("H25").Value = 1933
$ws.Range("K25").Value = 1933
$ws.Range("M25").Value = -1759
$ws.Range("I25").Value = 1933
$ws.Range("H35").Value = 5190.6665
$ws.Range("K35").Value = 5124.2
$ws.Range("M35").Value = -4830.2
$ws.Range("I35").Value = 5124.2
$ws.Range("H50").Value = 49631.668
$ws.Range("N50").Value = -50881.668
$ws.Range("J50").Value = 49631.668
$ws.Range("L50").Value = 49631.668
$ws.Range("H51").Value = 500
$ws.Range("K51").Value = 500
$ws.Range("M51").Value = 236
$ws.Range("I51").Value = 500
$ws.Range("H61").Value = 500
$ws.Range("K61").Value = 500
$ws.Range("M61").Value = -152
$ws.Range("I61").Value = 500
$ws.Range("H107").Value = 932.6667
$ws.Range("K107").Value = 800
$ws.Range("M107").Value = 1120
$ws.Range("I107").Value = 800

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2885.524
$ws.Range("N68").Value = -10905.8948
$ws.Range("J68").Value = 3094.6316
$ws.Range("L68").Value = 9283.8948
$ws.Range("H70").Value = 7999.4
$ws.Range("K70").Value = 23998.2
$ws.Range("M70").Value = -23683.2
$ws.Range("I70").Value = 7999.4
$ws.Range("H71").Value = 2885.524
$ws.Range("N71").Value = -35963.6844
$ws.Range("J71").Value = 3094.6316
$ws.Range("L71").Value = 27851.6844
$ws.Range("H73").Value = 7999.4
$ws.Range("K73").Value = 23998.2
$ws.Range("M73").Value = -22906.2
$ws.Range("I73").Value = 7999.4
$ws.Range("H97").Value = 1199.3334
$ws.Range("N97").Value = -6992
$ws.Range("K97").Value = 2397
$ws.Range("J97").Value = 2000
$ws.Range("M97").Value = -1901
$ws.Range("L97").Value = 6000
$ws.Range("I97").Value = 799
$ws.Range("H107").Value = 334.83334
$ws.Range("N107").Value = -4815
$ws.Range("K107").Value = 1010.4
$ws.Range("J107").Value = 325
$ws.Range("M107").Value = 909.5999999999999
$ws.Range("L107").Value = 975
$ws.Range("I107").Value = 336.8

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 67.75
$ws.Range("K2").Value = 70.8
$ws.Range("M2").Value = 42.2
$ws.Range("I2").Value = 70.8
$ws.Range("H97").Value = 2000
$ws.Range("N97").Value = -2992
$ws.Range("J97").Value = 2000
$ws.Range("L97").Value = 2000
$ws.Range("H132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1500
$ws.Range("K7").Value = 1500
$ws.Range("M7").Value = -1388
$ws.Range("I7").Value = 1500
$ws.Range("H46").Value = 251130.25
$ws.Range("N46").Value = -2622
$ws.Range("K46").Value = 334091.66
$ws.Range("J46").Value = 2246
$ws.Range("M46").Value = -333903.66
$ws.Range("L46").Value = 2246
$ws.Range("I46").Value = 334091.66
$ws.Range("H82").Value = 1412.7368
$ws.Range("K82").Value = 1367.1765
$ws.Range("M82").Value = -1006.1765
$ws.Range("I82").Value = 1367.1765
$ws.Range("H85").Value = 1412.7368
$ws.Range("K85").Value = 1367.1765
$ws.Range("M85").Value = -119.1765
$ws.Range("I85").Value = 1367.1765
$ws.Range("H122").Value = 4933.5
$ws.Range("N122").Value = -19883.5
$ws.Range("K122").Value = 14617.5
$ws.Range("J122").Value = 4994.5
$ws.Range("M122").Value = -12167.5
$ws.Range("L122").Value = 14983.5
$ws.Range("I122").Value = 4872.5
$ws.Range("H126").Value = 1500
$ws.Range("K126").Value = 4500
$ws.Range("M126").Value = -2030
$ws.Range("I126").Value = 1500
$ws.Range("H139").Value = 0
$ws.Range("N139").ClearContents()
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 25000
$ws.Range("K32").Value = 25000
$ws.Range("M32").Value = -24683
$ws.Range("I32").Value = 25000
$ws.Range("H95").Value = 50000
$ws.Range("N95").Value = -55492
$ws.Range("J95").Value = 50000
$ws.Range("L95").Value = 50000
$ws.Range("H97").Value = 45666.668
$ws.Range("N97").Value = -47648.668
$ws.Range("J97").Value = 45666.668
$ws.Range("L97").Value = 45666.668
$ws.Range("H107").Value = 656
$ws.Range("N107").ClearContents()
$ws.Range("K107").Value = 1968
$ws.Range("J107").Value = 0
$ws.Range("M107").Value = -48
$ws.Range("L107").Value = 0
$ws.Range("I107").Value = 656
$ws.Range("H113").Value = 8457.714
$ws.Range("N113").Value = -33742.001
$ws.Range("J113").Value = 9800.666999999999
$ws.Range("L113").Value = 29402.001
$ws.Range("H129").Value = 0
$ws.Range("N129").ClearContents()
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("H132").Value = 1737.875
$ws.Range("N132").ClearContents()
$ws.Range("K132").Value = 5213.625
$ws.Range("J132").Value = 0
$ws.Range("M132").Value = -2683.625
$ws.Range("L132").Value = 0
$ws.Range("I132").Value = 1737.875

Write-Host "Applied all Golem_Profits updates"